$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for the three added columns (M, N, O)
$ws.Range("M1").Value = "vissystemsGT"
$ws.Range("N1").Value = "sensorsGT"
$ws.Range("O1").Value = "aiGT"

# New trend data pulled in from Google Trends ("z trendsów wyniki")
$data = @(
    @(340, 507, 582),
    @(279, 487, 639),
    @(258, 455, 692),
    @(225, 419, 599),
    @(196, 407, 677),
    @(166, 414, 601),
    @(153, 432, 659),
    @(138, 455, 675),
    @(127, 455, 635),
    @(127, 478, 795),
    @(106, 478, 888),
    @(102, 475, 952)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 13).Value = $data[$i][0]
    $ws.Cells.Item($row, 14).Value = $data[$i][1]
    $ws.Cells.Item($row, 15).Value = $data[$i][2]
}

$ws.Range("L15").Select() | Out-Null
